# Fill in the previously-empty row 53 of the "Journal de travail" sheet
# with a new journal entry, matching the author's commit:
#   feat: created a custom error handler

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Date (2023-06-12), Type, Temps [h] and the long "Travail effectué" text
$ws.Range("A53").Value = 45089
$ws.Range("B53").Value = "Implémentation"
$ws.Range("C53").Value = 6

$ws.Range("D53").Value = "Frontend: Correction authentification, ajout i18n, fixes & refactor`nBackend: meilleure gestion des erreurs, ajout champs actif sur écran"
$ws.Range("D53").WrapText = $true
$ws.Rows.Item(53).RowHeight = 34

# Move the view to where the new entry was added (best-effort; harmless if unsupported)
try {
    $excel.ActiveWindow.ScrollRow = 34
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("D56").Select() | Out-Null
